$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 104.794801
$ws.Range("H2").Value = 314.384403
$ws.Range("I2").Value = 0.3872421191355361
$ws.Range("J2").Value = 0.3872421191355361
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.01567066666667
$ws.Range("N2").Value = 75.047012
$ws.Range("O2").Value = 0.3530689998156723
$ws.Range("P2").Value = 0.3530689998156723
$ws.Range("Q2").Value = 2621.512229394871
$ws.Range("R2").Value = 23593.61006455384
$ws.Range("S2").Value = 0.1367231876896852
$ws.Range("T2").Value = 0.1367231876896852
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 104.794801
$ws.Range("H3").Value = 314.384403
$ws.Range("I3").Value = 0.3872421191355361
$ws.Range("J3").Value = 0.3872421191355361
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.99788533333333
$ws.Range("N3").Value = 77.993656
$ws.Range("O3").Value = 0.3669318921836303
$ws.Range("P3").Value = 0.3669318921836303
$ws.Range("Q3").Value = 2724.443219927486
$ws.Range("R3").Value = 24519.98897934737
$ws.Range("S3").Value = 0.1420914835076011
$ws.Range("T3").Value = 0.1420914835076011
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 104.794801
$ws.Range("H4").Value = 314.384403
$ws.Range("I4").Value = 0.3872421191355361
$ws.Range("J4").Value = 0.3872421191355361
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.83851733333333
$ws.Range("N4").Value = 59.515552
$ws.Range("O4").Value = 0.2799991080006974
$ws.Range("P4").Value = 0.2799991080006974
$ws.Range("Q4").Value = 2078.973476081717
$ws.Range("R4").Value = 18710.76128473546
$ws.Range("S4").Value = 0.1084274479382499
$ws.Range("T4").Value = 0.1084274479382499
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 136.674446
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5050451128841343
$ws.Range("J5").Value = 0.5050451128841343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.01567066666667
$ws.Range("N5").Value = 75.047012
$ws.Range("O5").Value = 0.3530689998156723
$ws.Range("P5").Value = 0.3530689998156723
$ws.Range("Q5").Value = 3419.002929685117
$ws.Range("R5").Value = 30771.02636716605
$ws.Range("S5").Value = 0.1783157728677946
$ws.Range("T5").Value = 0.1783157728677946
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 136.674446
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5050451128841343
$ws.Range("J6").Value = 0.5050451128841343
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 25.99788533333333
$ws.Range("N6").Value = 77.993656
$ws.Range("O6").Value = 0.3669318921836303
$ws.Range("P6").Value = 0.3669318921836303
$ws.Range("Q6").Value = 3553.246575104858
$ws.Range("R6").Value = 31979.21917594373
$ws.Range("S6").Value = 0.1853171589086706
$ws.Range("T6").Value = 0.1853171589086706
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 136.674446
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5050451128841343
$ws.Range("J7").Value = 0.5050451128841343
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.83851733333333
$ws.Range("N7").Value = 59.515552
$ws.Range("O7").Value = 0.2799991080006974
$ws.Range("P7").Value = 0.2799991080006974
$ws.Range("Q7").Value = 2711.41836599473
$ws.Range("R7").Value = 24402.76529395257
$ws.Range("S7").Value = 0.1414121811076691
$ws.Range("T7").Value = 0.1414121811076691
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.14904533333333
$ws.Range("H8").Value = 87.447136
$ws.Range("I8").Value = 0.1077127679803296
$ws.Range("J8").Value = 0.1077127679803296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.01567066666667
$ws.Range("N8").Value = 75.047012
$ws.Range("O8").Value = 0.3530689998156723
$ws.Range("P8").Value = 0.3530689998156723
$ws.Range("Q8").Value = 729.1829183064035
$ws.Range("R8").Value = 6562.646264757632
$ws.Range("S8").Value = 0.03803003925819253
$ws.Range("T8").Value = 0.03803003925819253
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.14904533333333
$ws.Range("H9").Value = 87.447136
$ws.Range("I9").Value = 0.1077127679803296
$ws.Range("J9").Value = 0.1077127679803296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.99788533333333
$ws.Range("N9").Value = 77.993656
$ws.Range("O9").Value = 0.3669318921836303
$ws.Range("P9").Value = 0.3669318921836303
$ws.Range("Q9").Value = 757.813538152135
$ws.Range("R9").Value = 6820.321843369216
$ws.Range("S9").Value = 0.03952324976735868
$ws.Range("T9").Value = 0.03952324976735868
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.14904533333333
$ws.Range("H10").Value = 87.447136
$ws.Range("I10").Value = 0.1077127679803296
$ws.Range("J10").Value = 0.1077127679803296
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.83851733333333
$ws.Range("N10").Value = 59.515552
$ws.Range("O10").Value = 0.2799991080006974
$ws.Range("P10").Value = 0.2799991080006974
$ws.Range("Q10").Value = 578.2738410954524
$ws.Range("R10").Value = 5204.464569859072
$ws.Range("S10").Value = 0.03015947895477836
$ws.Range("T10").Value = 0.03015947895477836
